$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 86 (Item ID 12603)
$ws.Range("H86").Value = 100641.4
$ws.Range("I86").Value = 50501.5
$ws.Range("J86").Value = 134068
$ws.Range("K86").Value = 50501.5
$ws.Range("L86").Value = 134068
$ws.Range("M86").Value = -49378.5
$ws.Range("N86").Value = -136314
# Row 89 (Item ID 12603)
$ws.Range("H89").Value = 100641.4
$ws.Range("I89").Value = 50501.5
$ws.Range("J89").Value = 134068
$ws.Range("K89").Value = 252507.5
$ws.Range("L89").Value = 670340
$ws.Range("M89").Value = -246891.5
$ws.Range("N89").Value = -681572
# Row 96 (Item ID 19894)
$ws.Range("H96").Value = 50002616
$ws.Range("I96").Value = 83336664
$ws.Range("J96").Value = 1540
$ws.Range("K96").Value = 250009992
$ws.Range("L96").Value = 4620
$ws.Range("M96").Value = -250008619
$ws.Range("N96").Value = -7366
# Row 100 (Item ID 19906)
$ws.Range("H100").Value = 37040116
$ws.Range("I100").Value = 66669012
$ws.Range("J100").Value = 4000
$ws.Range("K100").Value = 66669012
$ws.Range("L100").Value = 4000
$ws.Range("M100").Value = -66668471
$ws.Range("N100").Value = -5082
# Row 132 (Item ID 44049)
$ws.Range("H132").Value = 51163.65
$ws.Range("I132").Value = 56765.168
$ws.Range("J132").Value = 750
$ws.Range("K132").Value = 170295.504
$ws.Range("L132").Value = 2250
$ws.Range("M132").Value = -167765.504
$ws.Range("N132").Value = -7310
# Row 135 (Item ID 44047)
$ws.Range("H135").Value = 1262.25
$ws.Range("I135").Value = 854.2222
$ws.Range("J135").Value = 1786.8572
$ws.Range("K135").Value = 7687.999800000001
$ws.Range("L135").Value = 16081.7148
$ws.Range("M135").Value = -5152.999800000001
$ws.Range("N135").Value = -21151.7148
# Row 138 (Item ID 44169)
$ws.Range("H138").Value = 1677.3549
$ws.Range("I138").Value = 1273.7838
$ws.Range("J138").Value = 2274.64
$ws.Range("K138").Value = 3821.3514
$ws.Range("L138").Value = 6823.92
$ws.Range("M138").Value = 1318.6486
$ws.Range("N138").Value = -17103.92

$ws = $wb.Worksheets.Item("ARM")
# Row 74 (Item ID 44000)
$ws.Range("H74").Value = 1923.4286
$ws.Range("I74").Value = 2436.4443
$ws.Range("J74").Value = 1000
$ws.Range("K74").Value = 2436.4443
$ws.Range("L74").Value = 1000
$ws.Range("M74").Value = -1562.4443
$ws.Range("N74").Value = -2748
# Row 77 (Item ID 44000)
$ws.Range("H77").Value = 1923.4286
$ws.Range("I77").Value = 2436.4443
$ws.Range("J77").Value = 1000
$ws.Range("K77").Value = 12182.2215
$ws.Range("L77").Value = 5000
$ws.Range("M77").Value = -7814.2215
$ws.Range("N77").Value = -13736
# Row 92 (Item ID 18050)
$ws.Range("H92").Value = 17687.75
$ws.Range("J92").Value = 17687.75
$ws.Range("L92").Value = 17687.75
$ws.Range("N92").Value = -22679.75
# Row 97 (Item ID 19941)
$ws.Range("H97").Value = 1637.1786
$ws.Range("I97").Value = 1281.2609
$ws.Range("J97").Value = 3274.4
$ws.Range("K97").Value = 1281.2609
$ws.Range("L97").Value = 3274.4
$ws.Range("M97").Value = -785.2609
$ws.Range("N97").Value = -4266.4
# Row 122 (Item ID 36168)
$ws.Range("H122").Value = 2571.5386
$ws.Range("I122").Value = 2118.5715
$ws.Range("J122").Value = 3100
$ws.Range("K122").Value = 6355.7145
$ws.Range("L122").Value = 9300
$ws.Range("M122").Value = -3905.7145
$ws.Range("N122").Value = -14200

$ws = $wb.Worksheets.Item("BSM")
# Row 86 (Item ID 12526)
$ws.Range("H86").Value = 3500
$ws.Range("I86").Value = 5000
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 5000
$ws.Range("L86").Value = 2000
$ws.Range("M86").Value = -3877
$ws.Range("N86").Value = -4246
# Row 89 (Item ID 12526)
$ws.Range("H89").Value = 3500
$ws.Range("I89").Value = 5000
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 25000
$ws.Range("L89").Value = 10000
$ws.Range("M89").Value = -19384
$ws.Range("N89").Value = -21232
# Row 105 (Item ID 19947)
$ws.Range("H105").Value = 3092
$ws.Range("I105").Value = 2869.1052
$ws.Range("J105").Value = 4503.6665
$ws.Range("K105").Value = 2869.1052
$ws.Range("L105").Value = 4503.6665
$ws.Range("M105").Value = -1122.1052
$ws.Range("N105").Value = -7997.6665

$ws = $wb.Worksheets.Item("CRP")
# Row 58 (Item ID 44021)
$ws.Range("H58").Value = 1979.9395
$ws.Range("I58").Value = 1701.8889
$ws.Range("J58").Value = 2313.6
$ws.Range("K58").Value = 1701.8889
$ws.Range("L58").Value = 2313.6
$ws.Range("M58").Value = -1498.8889
$ws.Range("N58").Value = -2719.6
# Row 74 (Item ID 10636)
$ws.Range("H74").Value = 38500
$ws.Range("J74").Value = 38500
$ws.Range("L74").Value = 38500
$ws.Range("N74").Value = -40248
# Row 77 (Item ID 10636)
$ws.Range("H77").Value = 38500
$ws.Range("J77").Value = 38500
$ws.Range("L77").Value = 115500
$ws.Range("N77").Value = -124236
# Row 81 (Item ID 10613)
$ws.Range("H81").Value = 28000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 28000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 28000
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -29996
# Row 84 (Item ID 10613)
$ws.Range("H84").Value = 28000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 28000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 84000
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -93984
# Row 105 (Item ID 19928)
$ws.Range("H105").Value = 1187
$ws.Range("I105").Value = 1187
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1187
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 560
$ws.Range("N105").ClearContents()
# Row 109 (Item ID 27203)
$ws.Range("H109").Value = 21183.334
$ws.Range("J109").Value = 21183.334
$ws.Range("L109").Value = 21183.334
$ws.Range("N109").Value = -23263.334
# Row 132 (Item ID 44019)
$ws.Range("H132").Value = 2286.9614
$ws.Range("I132").Value = 1092.3158
$ws.Range("J132").Value = 5529.5713
$ws.Range("K132").Value = 3276.9474
$ws.Range("L132").Value = 16588.7139
$ws.Range("M132").Value = -746.9474
$ws.Range("N132").Value = -21648.7139
# Row 136 (Item ID 44021)
$ws.Range("H136").Value = 1979.9395
$ws.Range("I136").Value = 1701.8889
$ws.Range("J136").Value = 2313.6
$ws.Range("K136").Value = 5105.6667
$ws.Range("L136").Value = 6940.799999999999
$ws.Range("M136").Value = -2555.6667
$ws.Range("N136").Value = -12040.8

$ws = $wb.Worksheets.Item("CUL")
# Row 131 (Item ID 36060)
$ws.Range("H131").Value = 5112.593
$ws.Range("J131").Value = 5490
$ws.Range("L131").Value = 16470
$ws.Range("N131").Value = -26550
# Row 133 (Item ID 44073)
$ws.Range("H133").Value = 2350.7693
$ws.Range("I133").Value = 2912
$ws.Range("K133").Value = 8736
$ws.Range("M133").Value = -3676
# Row 134 (Item ID 44074)
$ws.Range("H134").Value = 2096
$ws.Range("I134").Value = 2320
$ws.Range("K134").Value = 6960
$ws.Range("M134").Value = -1890

$ws = $wb.Worksheets.Item("GSM")
# Row 41 (Item ID 2449)
$ws.Range("H41").Value = 1207
$ws.Range("I41").Value = 560.5
$ws.Range("J41").Value = 2500
$ws.Range("K41").Value = 560.5
$ws.Range("L41").Value = 2500
$ws.Range("M41").Value = -205.5
$ws.Range("N41").Value = -3210
# Row 132 (Item ID 44008)
$ws.Range("H132").Value = 2358.4
$ws.Range("I132").Value = 1997
$ws.Range("K132").Value = 5991
$ws.Range("M132").Value = -3461

$ws = $wb.Worksheets.Item("LTW")
# Row 46 (Item ID 5282)
$ws.Range("H46").Value = 32507.125
$ws.Range("I46").Value = 72475.86
$ws.Range("J46").Value = 1420.3334
$ws.Range("K46").Value = 72475.86
$ws.Range("L46").Value = 1420.3334
$ws.Range("M46").Value = -72287.86
$ws.Range("N46").Value = -1796.3334
# Row 93 (Item ID 19993)
$ws.Range("H93").Value = 3879.9033
$ws.Range("I93").Value = 4542.76
$ws.Range("J93").Value = 1118
$ws.Range("K93").Value = 4542.76
$ws.Range("L93").Value = 1118
$ws.Range("M93").Value = -3294.76
$ws.Range("N93").Value = -3614

$ws = $wb.Worksheets.Item("WVR")
# Row 132 (Item ID 44029)
$ws.Range("H132").Value = 3112.875
$ws.Range("I132").Value = 2918.9092
$ws.Range("J132").Value = 3539.6
$ws.Range("K132").Value = 8756.7276
$ws.Range("L132").Value = 10618.8
$ws.Range("M132").Value = -6226.7276
$ws.Range("N132").Value = -15678.8
# Row 136 (Item ID 44031)
$ws.Range("H136").Value = 3281.0833
$ws.Range("J136").Value = 4931.077
$ws.Range("L136").Value = 14793.231
$ws.Range("N136").Value = -19893.231
